$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Troupes")

# Widen column A slightly (raw OOXML width 14.28515625 -> 16)
$ws.Columns.Item(1).ColumnWidth = 15.15

# Update the "19/joueurs*3" label and formulas to use 37 instead of 19
$ws.Range("A4").Value = "37/joueurs*3"
$ws.Range("B4").Formula = "=37/B1*3"
$ws.Range("C4").Formula = "=37/C1*3"
$ws.Range("D4").Formula = "=37/D1*3"

# Add a new row 5 (copy formatting of row 4) with the new objective line
$ws.Range("A4:D4").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A5").Value = "0.5*37/joueurs*4"
$ws.Range("B5").Formula = "=0.5*B4"
$ws.Range("C5:D5").Formula = "=0.5*C4"

# Move the active selection like in the saved file
$ws.Range("D9").Select()
